# Commit message: "Caculate C, D, T and plot data"
#
# This updates the "Time - uS" (H) column on the CaculateResponseTime sheet
# with recalculated response-time values, and nudges the active selection
# to where the user's cursor ended up after the edit.
#
# It also tidies the row heights on the "Hard" sheet (rows 1,3,4,5,7,8)
# back to their natural/auto height.

$wb = $excel.ActiveWorkbook

# --- CaculateResponseTime sheet: recalculated response times (column H) ---
$ws = $wb.Worksheets.Item("CaculateResponseTime")

$ws.Range("H2").Value  = 480
$ws.Range("H3").Value  = 494
$ws.Range("H4").Value  = 494
$ws.Range("H5").Value  = 494
$ws.Range("H6").Value  = 494
$ws.Range("H7").Value  = 494
$ws.Range("H8").Value  = 494
$ws.Range("H9").Value  = 494
$ws.Range("H10").Value = 494
$ws.Range("H11").Value = 494
$ws.Range("H12").Value = 494
$ws.Range("H13").Value = 499
$ws.Range("H14").Value = 735

# Leave the cursor where the author left it (J14) after finishing the edit.
$ws.Range("J14").Select() | Out-Null

# --- Hard sheet: row heights settle back to the natural (auto) height ---
$hard = $wb.Worksheets.Item("Hard")
$hard.Rows.Item(1).RowHeight = 13.2
$hard.Rows.Item(3).RowHeight = 13.2
$hard.Rows.Item(4).RowHeight = 13.2
$hard.Rows.Item(5).RowHeight = 13.2
$hard.Rows.Item(7).RowHeight = 13.2
$hard.Rows.Item(8).RowHeight = 13.2
